{"js": "// Update the worksheet date and all the two-digit\u00f7one-digit division\n// problems to the new values for the \"generated at c986bee\" refresh.\nconst replacements = [\n  [\"2024-09-16 Monday\", \"2024-09-17 Tuesday\"],\n  [\"34\u00f77=\", \"21\u00f77=\"],\n  [\"63\u00f74=\", \"66\u00f78=\"],\n  [\"93\u00f77=\", \"86\u00f79=\"],\n  [\"92\u00f79=\", \"28\u00f76=\"],\n  [\"63\u00f76=\", \"26\u00f77=\"],\n  [\"83\u00f74=\", \"29\u00f74=\"],\n  [\"37\u00f72=\", \"48\u00f73=\"],\n  [\"45\u00f75=\", \"36\u00f77=\"],\n  [\"34\u00f72=\", \"16\u00f74=\"],\n  [\"35\u00f78=\", \"44\u00f74=\"],\n  [\"72\u00f79=\", \"11\u00f76=\"],\n  [\"97\u00f73=\", \"66\u00f72=\"],\n  [\"24\u00f77=\", \"41\u00f73=\"],\n  [\"47\u00f77=\", \"27\u00f75=\"],\n  [\"42\u00f77=\", \"45\u00f72=\"],\n  [\"77\u00f79=\", \"89\u00f77=\"],\n  [\"52\u00f73=\", \"16\u00f75=\"],\n  [\"81\u00f76=\", \"41\u00f75=\"],\n  [\"69\u00f78=\", \"58\u00f78=\"],\n  [\"26\u00f73=\", \"52\u00f77=\"],\n  [\"37\u00f76=\", \"76\u00f76=\"],\n  [\"10\u00f74=\", \"72\u00f77=\"],\n  [\"33\u00f75=\", \"53\u00f79=\"],\n  [\"98\u00f77=\", \"93\u00f78=\"],\n  [\"85\u00f76=\", \"88\u00f73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and all the two-digit\u00f7one-digit division\n# problems to the new values for the \"generated at c986bee\" refresh.\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$replacements = @(\n  @{old=\"2024-09-16 Monday\"; new=\"2024-09-17 Tuesday\"},\n  @{old=\"34\u00f77=\"; new=\"21\u00f77=\"},\n  @{old=\"63\u00f74=\"; new=\"66\u00f78=\"},\n  @{old=\"93\u00f77=\"; new=\"86\u00f79=\"},\n  @{old=\"92\u00f79=\"; new=\"28\u00f76=\"},\n  @{old=\"63\u00f76=\"; new=\"26\u00f77=\"},\n  @{old=\"83\u00f74=\"; new=\"29\u00f74=\"},\n  @{old=\"37\u00f72=\"; new=\"48\u00f73=\"},\n  @{old=\"45\u00f75=\"; new=\"36\u00f77=\"},\n  @{old=\"34\u00f72=\"; new=\"16\u00f74=\"},\n  @{old=\"35\u00f78=\"; new=\"44\u00f74=\"},\n  @{old=\"72\u00f79=\"; new=\"11\u00f76=\"},\n  @{old=\"97\u00f73=\"; new=\"66\u00f72=\"},\n  @{old=\"24\u00f77=\"; new=\"41\u00f73=\"},\n  @{old=\"47\u00f77=\"; new=\"27\u00f75=\"},\n  @{old=\"42\u00f77=\"; new=\"45\u00f72=\"},\n  @{old=\"77\u00f79=\"; new=\"89\u00f77=\"},\n  @{old=\"52\u00f73=\"; new=\"16\u00f75=\"},\n  @{old=\"81\u00f76=\"; new=\"41\u00f75=\"},\n  @{old=\"69\u00f78=\"; new=\"58\u00f78=\"},\n  @{old=\"26\u00f73=\"; new=\"52\u00f77=\"},\n  @{old=\"37\u00f76=\"; new=\"76\u00f76=\"},\n  @{old=\"10\u00f74=\"; new=\"72\u00f77=\"},\n  @{old=\"33\u00f75=\"; new=\"53\u00f79=\"},\n  @{old=\"98\u00f77=\"; new=\"93\u00f78=\"},\n  @{old=\"85\u00f76=\"; new=\"88\u00f73=\"}\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($pair.old, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $pair.new, $wdReplaceAll) | Out-Null\n}\n"}
